$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 63, shifting existing rows 63:70 down to 64:71
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly record
$ws.Cells.Item(63, 1).Value = 11
$ws.Cells.Item(63, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(63, 3).Value = "Bíobío"
$ws.Cells.Item(63, 4).Value = 44505
$ws.Cells.Item(63, 5).Value = 8
$ws.Cells.Item(63, 6).Value = 100112032
$ws.Cells.Item(63, 7).Value = "Zapallo italiano"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 350
$ws.Cells.Item(63, 11).Value = 6500
$ws.Cells.Item(63, 12).Value = 7000
$ws.Cells.Item(63, 13).Value = 6714
$ws.Cells.Item(63, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(63, 15).Value = "Región del Maule"
$ws.Cells.Item(63, 16).Value = 112
$ws.Cells.Item(63, 17).Value = 60
$ws.Cells.Item(63, 18).Value = "Hortaliza"
